$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 60 (pushes the existing rows 60-92 down to 61-93),
# then populate the new row with the "Alcachofa" / "Argentina(o)" entry.
$ws.Rows.Item(60).Insert()

$ws.Range("A60").Value2 = 5
$ws.Range("B60").Value2 = "Macroferia Regional de Talca"
$ws.Range("C60").Value2 = "Maule"
$ws.Range("D60").Value2 = 44784
$ws.Range("E60").Value2 = 7
$ws.Range("F60").Value2 = 100112013
$ws.Range("G60").Value2 = "Alcachofa"
$ws.Range("H60").Value2 = "Argentina(o)"
$ws.Range("I60").Value2 = "Primera"
$ws.Range("J60").Value2 = 300
$ws.Range("K60").Value2 = 13000
$ws.Range("L60").Value2 = 13000
$ws.Range("M60").Value2 = 13000
$ws.Range("N60").Value2 = "$/caja 40 unidades"
$ws.Range("O60").Value2 = "Provincia del Elquí"
$ws.Range("P60").Value2 = 325
$ws.Range("Q60").Value2 = 40
$ws.Range("R60").Value2 = "Hortaliza"
